$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.268.20"
Set-TextValue $ws.Range("E2") "  +1.27%  "
Set-TextValue $ws.Range("D3") "2.636.49"
Set-TextValue $ws.Range("E3") "  +0.99%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "599.21"
Set-TextValue $ws.Range("E5") "  +1.24%  "
Set-TextValue $ws.Range("D6") "154.24"
Set-TextValue $ws.Range("E6") "  +2.34%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("E8") "  -0.12%  "
Set-TextValue $ws.Range("D9") "2.635.34"
Set-TextValue $ws.Range("E9") "  +0.98%  "
Set-TextValue $ws.Range("E10") "  +7.62%  "
Set-TextValue $ws.Range("E11") "  -0.65%  "
Set-TextValue $ws.Range("D12") "5.23"
Set-TextValue $ws.Range("E12") "  +1.04%  "
Set-TextValue $ws.Range("E13") "  +1.44%  "
Set-TextValue $ws.Range("D14") "27.97"
Set-TextValue $ws.Range("E14") "  +2.41%  "
Set-TextValue $ws.Range("E15") "  +2.90%  "
Set-TextValue $ws.Range("D17") "68.285.95"
Set-TextValue $ws.Range("E17") "  +1.45%  "
Set-TextValue $ws.Range("D18") "2.642.00"
Set-TextValue $ws.Range("D19") "11.40"
Set-TextValue $ws.Range("E19") "  +3.27%  "
Set-TextValue $ws.Range("D20") "366.00"
Set-TextValue $ws.Range("E20") "  -1.54%  "
Set-TextValue $ws.Range("D21") "7.41"
Set-TextValue $ws.Range("E21") "  +0.58%  "
Set-TextValue $ws.Range("E22") "  -0.54%  "
Set-TextValue $ws.Range("E23") "  +1.85%  "
Set-TextValue $ws.Range("E24") "  +3.69%  "
Set-TextValue $ws.Range("D25") "73.66"
Set-TextValue $ws.Range("E25") "  -0.03%  "
Set-TextValue $ws.Range("E26") "  -0.14%  "
Set-TextValue $ws.Range("D27") "10.05"
Set-TextValue $ws.Range("E27") "  +1.43%  "
Set-TextValue $ws.Range("D28") "2.770.28"
Set-TextValue $ws.Range("E29") "  +5.60%  "
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.19%  "
Set-TextValue $ws.Range("D31") "574.14"
Set-TextValue $ws.Range("E31") "  -1.07%  "
Set-TextValue $ws.Range("B32") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D32") "8.00"
Set-TextValue $ws.Range("E32") "  +4.26%  "
Set-TextValue $ws.Range("B33") "Fetch.AI"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D33") "1.41"
Set-TextValue $ws.Range("E33") "  +3.97%  "
Set-TextValue $ws.Range("E34") "  +2.48%  "
Set-TextValue $ws.Range("D35") "0.130"
Set-TextValue $ws.Range("E35") "  +2.72%  "
Set-TextValue $ws.Range("D36") "1.00"
Set-TextValue $ws.Range("E36") "  +0.00%  "
Set-TextValue $ws.Range("E37") "  +3.25%  "
Set-TextValue $ws.Range("D38") "160.33"
Set-TextValue $ws.Range("E38") "  +1.56%  "
Set-TextValue $ws.Range("D39") "19.26"
Set-TextValue $ws.Range("E39") "  +1.01%  "
Set-TextValue $ws.Range("E40") "  +2.95%  "
Set-TextValue $ws.Range("E41") "  +0.72%  "
Set-TextValue $ws.Range("E42") "  +2.74%  "
Set-TextValue $ws.Range("D43") "17.74"
Set-TextValue $ws.Range("E44") "  +2.75%  "
Set-TextValue $ws.Range("D45") "0.0₆0322"
Set-TextValue $ws.Range("E45") "  +14.09%  "
Set-TextValue $ws.Range("E46") "  -0.02%  "
Set-TextValue $ws.Range("D47") "40.50"
Set-TextValue $ws.Range("E47") "  -0.40%  "
Set-TextValue $ws.Range("D48") "157.02"
Set-TextValue $ws.Range("E48") "  +2.37%  "
Set-TextValue $ws.Range("D49") "3.75"
Set-TextValue $ws.Range("E49") "  +1.01%  "
Set-TextValue $ws.Range("D50") "1.71"
Set-TextValue $ws.Range("E50") "  +1.80%  "
Set-TextValue $ws.Range("D51") "21.86"
Set-TextValue $ws.Range("E51") "  +2.32%  "
